$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference date-format cell to copy number format from for new rows
$dateFormat = $ws.Cells.Item(2, 4).NumberFormat

# Row 54: Calidad=Especial
$ws.Cells.Item(54, 1).Value = 8
$ws.Cells.Item(54, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(54, 3).Value = "Coquimbo"
$ws.Cells.Item(54, 4).Value = 44476
$ws.Cells.Item(54, 5).Value = 4
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100107
$ws.Cells.Item(54, 8).Value = "Otros"
$ws.Cells.Item(54, 9).Value = 100107002
$ws.Cells.Item(54, 10).Value = "Chirimoya"
$ws.Cells.Item(54, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(54, 12).Value = "Especial"
$ws.Cells.Item(54, 13).Value = 240
$ws.Cells.Item(54, 14).Value = 2200
$ws.Cells.Item(54, 15).Value = 2300
$ws.Cells.Item(54, 16).Value = 2250
$ws.Cells.Item(54, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(54, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(54, 19).Value = 2250
$ws.Cells.Item(54, 20).Value = 1

# Row 55: Calidad=Primera
$ws.Cells.Item(55, 1).Value = 8
$ws.Cells.Item(55, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(55, 3).Value = "Coquimbo"
$ws.Cells.Item(55, 4).Value = 44476
$ws.Cells.Item(55, 5).Value = 4
$ws.Cells.Item(55, 6).Value = "Fruta"
$ws.Cells.Item(55, 7).Value = 100107
$ws.Cells.Item(55, 8).Value = "Otros"
$ws.Cells.Item(55, 9).Value = 100107002
$ws.Cells.Item(55, 10).Value = "Chirimoya"
$ws.Cells.Item(55, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(55, 12).Value = "Primera"
$ws.Cells.Item(55, 13).Value = 300
$ws.Cells.Item(55, 14).Value = 1900
$ws.Cells.Item(55, 15).Value = 2000
$ws.Cells.Item(55, 16).Value = 1950
$ws.Cells.Item(55, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(55, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(55, 19).Value = 1950
$ws.Cells.Item(55, 20).Value = 1

# Row 56: Calidad=Segunda
$ws.Cells.Item(56, 1).Value = 8
$ws.Cells.Item(56, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(56, 3).Value = "Coquimbo"
$ws.Cells.Item(56, 4).Value = 44476
$ws.Cells.Item(56, 5).Value = 4
$ws.Cells.Item(56, 6).Value = "Fruta"
$ws.Cells.Item(56, 7).Value = 100107
$ws.Cells.Item(56, 8).Value = "Otros"
$ws.Cells.Item(56, 9).Value = 100107002
$ws.Cells.Item(56, 10).Value = "Chirimoya"
$ws.Cells.Item(56, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(56, 12).Value = "Segunda"
$ws.Cells.Item(56, 13).Value = 300
$ws.Cells.Item(56, 14).Value = 1400
$ws.Cells.Item(56, 15).Value = 1500
$ws.Cells.Item(56, 16).Value = 1450
$ws.Cells.Item(56, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(56, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(56, 19).Value = 1450
$ws.Cells.Item(56, 20).Value = 1

# Row 57: Calidad=Especial
$ws.Cells.Item(57, 1).Value = 8
$ws.Cells.Item(57, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(57, 3).Value = "Coquimbo"
$ws.Cells.Item(57, 4).Value = 44172
$ws.Cells.Item(57, 5).Value = 4
$ws.Cells.Item(57, 6).Value = "Fruta"
$ws.Cells.Item(57, 7).Value = 100107
$ws.Cells.Item(57, 8).Value = "Otros"
$ws.Cells.Item(57, 9).Value = 100107002
$ws.Cells.Item(57, 10).Value = "Chirimoya"
$ws.Cells.Item(57, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(57, 12).Value = "Especial"
$ws.Cells.Item(57, 13).Value = 200
$ws.Cells.Item(57, 14).Value = 14000
$ws.Cells.Item(57, 15).Value = 14500
$ws.Cells.Item(57, 16).Value = 14250
$ws.Cells.Item(57, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(57, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(57, 19).Value = 1781
$ws.Cells.Item(57, 20).Value = 8

# Row 58: Calidad=Primera
$ws.Cells.Item(58, 1).Value = 8
$ws.Cells.Item(58, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(58, 3).Value = "Coquimbo"
$ws.Cells.Item(58, 4).Value = 44172
$ws.Cells.Item(58, 5).Value = 4
$ws.Cells.Item(58, 6).Value = "Fruta"
$ws.Cells.Item(58, 7).Value = 100107
$ws.Cells.Item(58, 8).Value = "Otros"
$ws.Cells.Item(58, 9).Value = 100107002
$ws.Cells.Item(58, 10).Value = "Chirimoya"
$ws.Cells.Item(58, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(58, 12).Value = "Primera"
$ws.Cells.Item(58, 13).Value = 200
$ws.Cells.Item(58, 14).Value = 12000
$ws.Cells.Item(58, 15).Value = 12500
$ws.Cells.Item(58, 16).Value = 12250
$ws.Cells.Item(58, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(58, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(58, 19).Value = 1531
$ws.Cells.Item(58, 20).Value = 8

# Row 59: Calidad=Segunda
$ws.Cells.Item(59, 1).Value = 8
$ws.Cells.Item(59, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 44172
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100107
$ws.Cells.Item(59, 8).Value = "Otros"
$ws.Cells.Item(59, 9).Value = 100107002
$ws.Cells.Item(59, 10).Value = "Chirimoya"
$ws.Cells.Item(59, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(59, 12).Value = "Segunda"
$ws.Cells.Item(59, 13).Value = 200
$ws.Cells.Item(59, 14).Value = 9500
$ws.Cells.Item(59, 15).Value = 10000
$ws.Cells.Item(59, 16).Value = 9750
$ws.Cells.Item(59, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(59, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(59, 19).Value = 1219
$ws.Cells.Item(59, 20).Value = 8

# Row 60: Calidad=Especial
$ws.Cells.Item(60, 1).Value = 8
$ws.Cells.Item(60, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(60, 3).Value = "Coquimbo"
$ws.Cells.Item(60, 4).Value = 44466
$ws.Cells.Item(60, 5).Value = 4
$ws.Cells.Item(60, 6).Value = "Fruta"
$ws.Cells.Item(60, 7).Value = 100107
$ws.Cells.Item(60, 8).Value = "Otros"
$ws.Cells.Item(60, 9).Value = 100107002
$ws.Cells.Item(60, 10).Value = "Chirimoya"
$ws.Cells.Item(60, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(60, 12).Value = "Especial"
$ws.Cells.Item(60, 13).Value = 160
$ws.Cells.Item(60, 14).Value = 2100
$ws.Cells.Item(60, 15).Value = 2200
$ws.Cells.Item(60, 16).Value = 2150
$ws.Cells.Item(60, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(60, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(60, 19).Value = 2150
$ws.Cells.Item(60, 20).Value = 1
$ws.Cells.Item(60, 4).NumberFormat = $dateFormat

# Row 61: Calidad=Primera
$ws.Cells.Item(61, 1).Value = 8
$ws.Cells.Item(61, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(61, 3).Value = "Coquimbo"
$ws.Cells.Item(61, 4).Value = 44466
$ws.Cells.Item(61, 5).Value = 4
$ws.Cells.Item(61, 6).Value = "Fruta"
$ws.Cells.Item(61, 7).Value = 100107
$ws.Cells.Item(61, 8).Value = "Otros"
$ws.Cells.Item(61, 9).Value = 100107002
$ws.Cells.Item(61, 10).Value = "Chirimoya"
$ws.Cells.Item(61, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(61, 12).Value = "Primera"
$ws.Cells.Item(61, 13).Value = 240
$ws.Cells.Item(61, 14).Value = 1700
$ws.Cells.Item(61, 15).Value = 1800
$ws.Cells.Item(61, 16).Value = 1750
$ws.Cells.Item(61, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(61, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(61, 19).Value = 1750
$ws.Cells.Item(61, 20).Value = 1
$ws.Cells.Item(61, 4).NumberFormat = $dateFormat

# Row 62: Calidad=Segunda
$ws.Cells.Item(62, 1).Value = 8
$ws.Cells.Item(62, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(62, 3).Value = "Coquimbo"
$ws.Cells.Item(62, 4).Value = 44466
$ws.Cells.Item(62, 5).Value = 4
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100107
$ws.Cells.Item(62, 8).Value = "Otros"
$ws.Cells.Item(62, 9).Value = 100107002
$ws.Cells.Item(62, 10).Value = "Chirimoya"
$ws.Cells.Item(62, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(62, 12).Value = "Segunda"
$ws.Cells.Item(62, 13).Value = 200
$ws.Cells.Item(62, 14).Value = 1300
$ws.Cells.Item(62, 15).Value = 1400
$ws.Cells.Item(62, 16).Value = 1350
$ws.Cells.Item(62, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(62, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(62, 19).Value = 1350
$ws.Cells.Item(62, 20).Value = 1
$ws.Cells.Item(62, 4).NumberFormat = $dateFormat

Write-Host "UsedRange now:" $ws.UsedRange.Address()
